$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3502.7778
$ws.Range("J17").Value = 2870.25
$ws.Range("L17").Value = 8610.75
$ws.Range("N17").Value = -8946.75
$ws.Range("H33").Value = 177.81818
$ws.Range("J33").Value = 75
$ws.Range("L33").Value = 75
$ws.Range("N33").Value = -533
$ws.Range("H96").Value = 1492.2
$ws.Range("I96").Value = 1290.1666
$ws.Range("J96").Value = 1626.8889
$ws.Range("K96").Value = 3870.4998
$ws.Range("L96").Value = 4880.6667
$ws.Range("M96").Value = -2497.4998
$ws.Range("N96").Value = -7626.6667
$ws.Range("H100").Value = 3395.2
$ws.Range("I100").Value = 998
$ws.Range("K100").Value = 998
$ws.Range("M100").Value = -457
$ws.Range("H111").Value = 2029
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4596.9346
$ws.Range("I32").Value = 3093.8572
$ws.Range("K32").Value = 3093.8572
$ws.Range("M32").Value = -2806.8572
$ws.Range("H45").Value = 6924592.5
$ws.Range("I45").Value = 22500752
$ws.Range("K45").Value = 22500752
$ws.Range("M45").Value = -22500375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 255038.75
$ws.Range("J86").Value = 402302
$ws.Range("L86").Value = 402302
$ws.Range("N86").Value = -404548
$ws.Range("H89").Value = 255038.75
$ws.Range("J89").Value = 402302
$ws.Range("L89").Value = 2011510
$ws.Range("N89").Value = -2022742
$ws.Range("H99").Value = 1577.5
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H138").Value = 29921.445
$ws.Range("J138").Value = 29921.445
$ws.Range("L138").Value = 29921.445
$ws.Range("N138").Value = -40201.445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 343.5
$ws.Range("I7").Value = 365.25
$ws.Range("K7").Value = 365.25
$ws.Range("M7").Value = -252.25
$ws.Range("H31").Value = 2317.647
$ws.Range("I31").Value = 2125
$ws.Range("K31").Value = 2125
$ws.Range("M31").Value = -1830
$ws.Range("H34").Value = 2317.647
$ws.Range("I34").Value = 2125
$ws.Range("K34").Value = 2125
$ws.Range("M34").Value = -1923
$ws.Range("H86").Value = 3277.5
$ws.Range("I86").Value = 3277.5
$ws.Range("K86").Value = 3277.5
$ws.Range("M86").Value = -2154.5
$ws.Range("H87").Value = 39000
$ws.Range("J87").Value = 39000
$ws.Range("L87").Value = 39000
$ws.Range("N87").Value = -41372
$ws.Range("H89").Value = 3277.5
$ws.Range("I89").Value = 3277.5
$ws.Range("K89").Value = 16387.5
$ws.Range("M89").Value = -10771.5
$ws.Range("H90").Value = 39000
$ws.Range("J90").Value = 39000
$ws.Range("L90").Value = 117000
$ws.Range("N90").Value = -128856
$ws.Range("H141").Value = 63799.8
$ws.Range("J141").Value = 61749.75
$ws.Range("L141").Value = 61749.75
$ws.Range("N141").Value = -72109.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 98.111115
$ws.Range("I12").Value = 58.25
$ws.Range("J12").Value = 130
$ws.Range("K12").Value = 174.75
$ws.Range("L12").Value = 390
$ws.Range("M12").Value = -1.75
$ws.Range("N12").Value = -736
$ws.Range("H40").Value = 130
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H131").Value = 12690.135
$ws.Range("J131").Value = 13458.396
$ws.Range("L131").Value = 40375.188
$ws.Range("N131").Value = -50455.188
$ws.Range("H133").Value = 4218.5713
$ws.Range("I133").Value = 2265
$ws.Range("K133").Value = 6795
$ws.Range("M133").Value = -1735

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 35000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 35000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 35000
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -35576
$ws.Range("H70").Value = 3966.6667
$ws.Range("I70").Value = 3800
$ws.Range("J70").Value = 4133.3335
$ws.Range("K70").Value = 3800
$ws.Range("L70").Value = 4133.3335
$ws.Range("M70").Value = -3530
$ws.Range("N70").Value = -4673.3335
$ws.Range("H73").Value = 3966.6667
$ws.Range("I73").Value = 3800
$ws.Range("J73").Value = 4133.3335
$ws.Range("K73").Value = 3800
$ws.Range("L73").Value = 4133.3335
$ws.Range("M73").Value = -2864
$ws.Range("N73").Value = -6005.3335
$ws.Range("H81").Value = 35000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 35000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 35000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -36996
$ws.Range("H84").Value = 35000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 35000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 105000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -114984
$ws.Range("H97").Value = 582.29034
$ws.Range("I97").Value = 576.7
$ws.Range("J97").Value = 750
$ws.Range("K97").Value = 576.7
$ws.Range("L97").Value = 750
$ws.Range("M97").Value = -80.70000000000005
$ws.Range("N97").Value = -1742
$ws.Range("H132").Value = 2265197
$ws.Range("I132").Value = 3206871
$ws.Range("J132").Value = 5179.2
$ws.Range("K132").Value = 9620613
$ws.Range("L132").Value = 15537.6
$ws.Range("M132").Value = -9618083
$ws.Range("N132").Value = -20597.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9931.474
$ws.Range("I40").Value = 9793.143
$ws.Range("J40").Value = 10318.8
$ws.Range("K40").Value = 9793.143
$ws.Range("L40").Value = 10318.8
$ws.Range("M40").Value = -9657.143
$ws.Range("N40").Value = -10590.8
$ws.Range("H46").Value = 1962.5834
$ws.Range("I46").Value = 1110.2858
$ws.Range("J46").Value = 3155.8
$ws.Range("K46").Value = 1110.2858
$ws.Range("L46").Value = 3155.8
$ws.Range("M46").Value = -922.2858000000001
$ws.Range("N46").Value = -3531.8
$ws.Range("H93").Value = 17544588
$ws.Range("I93").Value = 793.4
$ws.Range("J93").Value = 37037692
$ws.Range("K93").Value = 793.4
$ws.Range("L93").Value = 37037692
$ws.Range("M93").Value = 454.6
$ws.Range("N93").Value = -37040188
$ws.Range("H100").Value = 1894.6
$ws.Range("I100").Value = 896.5
$ws.Range("K100").Value = 896.5
$ws.Range("M100").Value = -355.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H70").Value = 29110
$ws.Range("J70").Value = 29110
$ws.Range("L70").Value = 29110
$ws.Range("N70").Value = -29740
$ws.Range("H73").Value = 29110
$ws.Range("J73").Value = 29110
$ws.Range("L73").Value = 29110
$ws.Range("N73").Value = -31294
$ws.Range("H81").Value = 1499
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 1499
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H132").Value = 1251.1333
$ws.Range("I132").Value = 903.61365
$ws.Range("K132").Value = 2710.84095
$ws.Range("M132").Value = -180.8409499999998
$ws.Range("H136").Value = 1749.8889
$ws.Range("I136").Value = 1436.4546
$ws.Range("J136").Value = 2242.4285
$ws.Range("K136").Value = 4309.3638
$ws.Range("L136").Value = 6727.2855
$ws.Range("M136").Value = -1759.3638
$ws.Range("N136").Value = -11827.2855
$ws.Range("H139").Value = 69846.75
$ws.Range("J139").Value = 69846.75
$ws.Range("L139").Value = 69846.75
$ws.Range("N139").Value = -80126.75
